# Author's edit: the "Planos:" label was changed to the singular "Plano:"
# (the rest of the document's wording is unchanged; the surrounding XML
# churn visible in the raw diff is just round-trip noise from the tool
# that produced the commit, not an intentional content/formatting change).

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "Planos: (",   # FindText
    $true,         # MatchCase
    $false,        # MatchWholeWord
    $false,        # MatchWildcards
    $false,        # MatchSoundsLike
    $false,        # MatchAllWordForms
    $true,         # Forward
    1,             # Wrap (wdFindContinue)
    $false,        # Format
    "Plano: (",    # ReplaceWith
    2              # Replace (wdReplaceAll)
)
